# Generate Report for Handback
#
# The row for "eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md" has moved from
# "Ready for handoff" to "Handed back: in sync with en-US" on every sheet,
# the zh-cn / de-de "Latest Handback DateTime" got refreshed, and the old
# "version mismatch" Error Detail message got cleared out now that the
# handback is up to date.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns for the eaeb81fe... row (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

# --- zh-cn sheet: Status / Latest Handback DateTime / Error Detail for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("K3").Value = "2016-09-04 16:53:43"
$zhcn.Range("P3").Value = ""

# --- de-de sheet: Status / Latest Handback DateTime / Error Detail for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $handedBack
$dede.Range("K3").Value = "2016-09-04 16:53:50"
$dede.Range("P3").Value = ""

# Error Detail column (P) is empty now, so it no longer needs to be as wide.
$zhcn.Columns.Item(16).ColumnWidth = 12.8
$dede.Columns.Item(16).ColumnWidth = 12.8
